# Applies the changes described by the commit diff:
#  - Shared string "newUser34422" -> "testSelenium29409" (Register!I2, LogIn!A2)
#  - Register sheet: scroll/selection change (topLeftCell=B1, activeCell/sqref=H16)
#  - LogIn sheet: selection change (activeCell/sqref=A5)

$wb = $excel.ActiveWorkbook

# --- Update the shared "newUser34422" value wherever it appears ---
$wsRegister = $wb.Worksheets.Item("Register")
$wsRegister.Range("I2").Value = "testSelenium29409"

$wsLogIn = $wb.Worksheets.Item("LogIn")
$wsLogIn.Range("A2").Value = "testSelenium29409"

# --- Update the view/selection on the Register sheet ---
$wsRegister.Select()
$wsRegister.Application.ActiveWindow.ScrollColumn = 2
$wsRegister.Range("H16").Select()

# --- Update the view/selection on the LogIn sheet ---
$wsLogIn.Select()
$wsLogIn.Range("A5").Select()
